$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data cells I2 and J2, both numeric value 8
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

# Match the header formatting (bold font, thin box border, centered/top
# aligned) already used by the other header cells (e.g. H1) by copying
# that cell's format onto the two new header cells.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
